# Supplemental_Tables.xlsx update — "updated markdown based on manuscript resubmission"
# Applies to "Table 1 - Caribbean warming" (the first worksheet):
#   - widen column A
#   - add a new "OISST" data column (G)
#   - convert the numeric measurement cells to text values
#   - add a new summary row for MHW frequency

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table 1 - Caribbean warming")

# --- widen column A (38.71 -> 55.71 stored width) -------------------------
$ws.Columns.Item(1).ColumnWidth = 54.83

# --- stage the whole table as Text-formatted so the string values we are  --
# --- about to write are stored as inline/shared strings (t="s") rather    --
# --- than being re-parsed back into numbers -------------------------------
$dataRange = $ws.Range("A1:G6")
$dataRange.NumberFormat = "@"

# --- header row -------------------------------------------------------------
$ws.Cells.Item(1,7).Value = "OISST"

# --- row 2: Caribbean Basin (°C per decade) ---------------------------------
$ws.Cells.Item(2,1).Value = "Caribbean Basin (°C per decade)"
$ws.Cells.Item(2,2).Value = "0.04"
$ws.Cells.Item(2,3).Value = "0.17"
$ws.Cells.Item(2,4).Value = "0.2"
$ws.Cells.Item(2,5).Value = "0.17"
$ws.Cells.Item(2,6).Value = "NA"
$ws.Cells.Item(2,7).Value = "NA"

# --- row 3: Caribbean Basin (total °C for period) ---------------------------
$ws.Cells.Item(3,1).Value = "Caribbean Basin (total °C for period)"
$ws.Cells.Item(3,2).Value = "0.6"
$ws.Cells.Item(3,3).Value = "0.68"
$ws.Cells.Item(3,4).Value = "0.54"
$ws.Cells.Item(3,5).Value = "0.66"
$ws.Cells.Item(3,6).Value = "NA"
$ws.Cells.Item(3,7).Value = "NA"

# --- row 4: Caribbean Reefs (°C per decade) ---------------------------------
$ws.Cells.Item(4,1).Value = "Caribbean Reefs (°C per decade)"
$ws.Cells.Item(4,2).Value = "0.04"
$ws.Cells.Item(4,3).Value = "0.15"
$ws.Cells.Item(4,4).Value = "0.17"
$ws.Cells.Item(4,5).Value = "0.19"
$ws.Cells.Item(4,6).Value = "0.18"
$ws.Cells.Item(4,7).Value = "NA"

# --- row 5: Caribbean Reefs (total °C for period) ---------------------------
$ws.Cells.Item(5,1).Value = "Caribbean Reefs (total °C for period)"
$ws.Cells.Item(5,2).Value = "0.6"
$ws.Cells.Item(5,3).Value = "0.6"
$ws.Cells.Item(5,4).Value = "0.46"
$ws.Cells.Item(5,5).Value = "0.74"
$ws.Cells.Item(5,6).Value = "0.47"
$ws.Cells.Item(5,7).Value = "NA"

# --- row 6 (new): Caribbean Basin (increasing frequency of MHW per year) ----
$ws.Cells.Item(6,1).Value = "Caribbean Basin (increasing frequency of MHW per year) "
$ws.Cells.Item(6,2).Value = "NA"
$ws.Cells.Item(6,3).Value = "NA"
$ws.Cells.Item(6,4).Value = "NA"
$ws.Cells.Item(6,5).Value = "NA"
$ws.Cells.Item(6,6).Value = "NA"
$ws.Cells.Item(6,7).Value = "0.05"

# --- drop the scratch "Text" number format again so cells keep the same  ---
# --- (default/general) style as the rest of the workbook -------------------
$dataRange.ClearFormats()
